# Apply scheduled market-price updates to the Odin_Profits workbook.
# Each block updates one leve row (H:N = price/profit columns) on its sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 39
$ws.Range("H39").Value = 581.3077
$ws.Range("I39").Value = 125.7
$ws.Range("J39").Value = 2100
$ws.Range("K39").Value = 377.1
$ws.Range("L39").Value = 6300
$ws.Range("M39").Value = -81.10000000000002
$ws.Range("N39").Value = -6892

# row 41
$ws.Range("H41").Value = 195.29411
$ws.Range("I41").Value = 258.75
$ws.Range("K41").Value = 258.75
$ws.Range("M41").Value = 181.25

# row 43
$ws.Range("H43").Value = 4835.222
$ws.Range("I43").Value = 5001
$ws.Range("J43").Value = 4814.5
$ws.Range("K43").Value = 5001
$ws.Range("L43").Value = 4814.5
$ws.Range("M43").Value = -4932
$ws.Range("N43").Value = -4952.5

# row 76
$ws.Range("H76").Value = 2317.7144
$ws.Range("I76").Value = 991.3333
$ws.Range("K76").Value = 991.3333
$ws.Range("M76").Value = -676.3333

# row 79
$ws.Range("H79").Value = 2317.7144
$ws.Range("I79").Value = 991.3333
$ws.Range("K79").Value = 991.3333
$ws.Range("M79").Value = 100.6667

# row 88
$ws.Range("H88").Value = 6499.5
$ws.Range("J88").Value = 6499.5
$ws.Range("L88").Value = 6499.5
$ws.Range("N88").Value = -7311.5

# row 91
$ws.Range("H91").Value = 6499.5
$ws.Range("J91").Value = 6499.5
$ws.Range("L91").Value = 6499.5
$ws.Range("N91").Value = -9307.5

# row 92
$ws.Range("H92").Value = 1002.619
$ws.Range("I92").Value = 897.6842
$ws.Range("K92").Value = 897.6842
$ws.Range("M92").Value = 350.3158

# row 100
$ws.Range("H100").Value = 3017.9583
$ws.Range("I100").Value = 956.7857
$ws.Range("K100").Value = 956.7857
$ws.Range("M100").Value = -415.7857

# row 116
$ws.Range("H116").Value = 10484.25
$ws.Range("I116").Value = 9601.237999999999
$ws.Range("K116").Value = 9601.237999999999
$ws.Range("M116").Value = -6159.237999999999

# row 135
$ws.Range("H135").Value = 2717.4707
$ws.Range("I135").Value = 1317.25
$ws.Range("K135").Value = 11855.25
$ws.Range("M135").Value = -9320.25

# row 138
$ws.Range("H138").Value = 3661.524
$ws.Range("I138").Value = 1013.7857
$ws.Range("J138").Value = 4985.393
$ws.Range("K138").Value = 3041.3571
$ws.Range("L138").Value = 14956.179
$ws.Range("M138").Value = 2098.6429
$ws.Range("N138").Value = -25236.179

$ws = $wb.Worksheets.Item("ARM")
# row 34
$ws.Range("H34").Value = 3000
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2729

# row 61
$ws.Range("H61").Value = 6647.8335
$ws.Range("J61").Value = 6647.8335
$ws.Range("L61").Value = 6647.8335
$ws.Range("N61").Value = -7071.8335

# row 74
$ws.Range("H74").Value = 2524.6667
$ws.Range("J74").Value = 2684
$ws.Range("L74").Value = 2684
$ws.Range("N74").Value = -4432

# row 77
$ws.Range("H77").Value = 2524.6667
$ws.Range("J77").Value = 2684
$ws.Range("L77").Value = 13420
$ws.Range("N77").Value = -22156

# row 97
$ws.Range("H97").Value = 1274.625
$ws.Range("I97").Value = 1140
$ws.Range("K97").Value = 1140
$ws.Range("M97").Value = -644

# row 132
$ws.Range("H132").Value = 70820.46000000001
$ws.Range("I132").Value = 2697.0908
$ws.Range("K132").Value = 8091.2724
$ws.Range("M132").Value = -5561.2724

# row 136
$ws.Range("H136").Value = 6647.8335
$ws.Range("J136").Value = 6647.8335
$ws.Range("L136").Value = 19943.5005
$ws.Range("N136").Value = -25043.5005

$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 4008.6428
$ws.Range("I86").Value = 2187.1
$ws.Range("K86").Value = 2187.1
$ws.Range("M86").Value = -1064.1

# row 89
$ws.Range("H89").Value = 4008.6428
$ws.Range("I89").Value = 2187.1
$ws.Range("K89").Value = 10935.5
$ws.Range("M89").Value = -5319.5

# row 99
$ws.Range("H99").Value = 7710.5815
$ws.Range("I99").Value = 7549.316
$ws.Range("J99").Value = 8936.200000000001
$ws.Range("K99").Value = 7549.316
$ws.Range("L99").Value = 8936.200000000001
$ws.Range("M99").Value = -6051.316
$ws.Range("N99").Value = -11932.2

# row 134
$ws.Range("H134").Value = 11430.538
$ws.Range("I134").Value = 5402.6665
$ws.Range("K134").Value = 16207.9995
$ws.Range("M134").Value = -13672.9995

$ws = $wb.Worksheets.Item("CRP")
# row 132
$ws.Range("H132").Value = 5312.65
$ws.Range("I132").Value = 4125
$ws.Range("J132").Value = 7094.125
$ws.Range("K132").Value = 12375
$ws.Range("L132").Value = 21282.375
$ws.Range("M132").Value = -9845
$ws.Range("N132").Value = -26342.375

$ws = $wb.Worksheets.Item("CUL")
# row 50
$ws.Range("H50").Value = 855.4286
$ws.Range("I50").Value = 759.6
$ws.Range("J50").Value = 1095
$ws.Range("K50").Value = 2278.8
$ws.Range("L50").Value = 3285
$ws.Range("M50").Value = -1797.8
$ws.Range("N50").Value = -4247

# row 53
$ws.Range("H53").Value = 855.4286
$ws.Range("I53").Value = 759.6
$ws.Range("J53").Value = 1095
$ws.Range("K53").Value = 2278.8
$ws.Range("L53").Value = 3285
$ws.Range("M53").Value = -1797.8
$ws.Range("N53").Value = -4247

# row 94
$ws.Range("H94").Value = 8225
$ws.Range("I94").Value = 3450
$ws.Range("J94").Value = 13000
$ws.Range("K94").Value = 10350
$ws.Range("L94").Value = 39000
$ws.Range("M94").Value = -9674
$ws.Range("N94").Value = -40352

# row 134
$ws.Range("H134").Value = 3656.3928
$ws.Range("I134").Value = 3382.4583
$ws.Range("K134").Value = 10147.3749
$ws.Range("M134").Value = -5077.374899999999

# row 140
$ws.Range("H140").Value = 23151012
$ws.Range("I140").Value = 28737980
$ws.Range("K140").Value = 86213940
$ws.Range("M140").Value = -86208760

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 7177.9585
$ws.Range("I80").Value = 4960.857
$ws.Range("J80").Value = 10281.9
$ws.Range("K80").Value = 4960.857
$ws.Range("L80").Value = 10281.9
$ws.Range("M80").Value = -3962.857
$ws.Range("N80").Value = -12277.9

# row 83
$ws.Range("H83").Value = 7177.9585
$ws.Range("I83").Value = 4960.857
$ws.Range("J83").Value = 10281.9
$ws.Range("K83").Value = 24804.285
$ws.Range("L83").Value = 51409.5
$ws.Range("M83").Value = -19812.285
$ws.Range("N83").Value = -61393.5

# row 102
$ws.Range("H102").Value = 6343.085
$ws.Range("I102").Value = 5439.893
$ws.Range("J102").Value = 7674.1055
$ws.Range("K102").Value = 5439.893
$ws.Range("L102").Value = 7674.1055
$ws.Range("M102").Value = -3817.893
$ws.Range("N102").Value = -10918.1055

# row 113
$ws.Range("H113").Value = 9310.846
$ws.Range("I113").Value = 3098
$ws.Range("J113").Value = 12072.111
$ws.Range("K113").Value = 3098
$ws.Range("L113").Value = 12072.111
$ws.Range("M113").Value = -928
$ws.Range("N113").Value = -16412.111

# row 126
$ws.Range("H126").Value = 7695.136
$ws.Range("I126").Value = 3092.8572
$ws.Range("K126").Value = 9278.571599999999
$ws.Range("M126").Value = -6808.571599999999

# row 132
$ws.Range("H132").Value = 4251.5415
$ws.Range("I132").Value = 3567.2354
$ws.Range("K132").Value = 10701.7062
$ws.Range("M132").Value = -8171.706200000001

$ws = $wb.Worksheets.Item("LTW")
# row 4
$ws.Range("H4").Value = 50000000
$ws.Range("I4").Value = 50000000
$ws.Range("K4").Value = 50000000
$ws.Range("M4").Value = -49999887

# row 28
$ws.Range("H28").Value = 50000000
$ws.Range("I28").Value = 50000000
$ws.Range("K28").Value = 50000000
$ws.Range("M28").Value = -49999768

# row 37
$ws.Range("H37").Value = 50000000
$ws.Range("I37").Value = 50000000
$ws.Range("K37").Value = 50000000
$ws.Range("M37").Value = -49999893

# row 82
$ws.Range("H82").Value = 5399.385
$ws.Range("I82").Value = 1064.1428
$ws.Range("J82").Value = 10457.167
$ws.Range("K82").Value = 1064.1428
$ws.Range("L82").Value = 10457.167
$ws.Range("M82").Value = -703.1428000000001
$ws.Range("N82").Value = -11179.167

# row 85
$ws.Range("H85").Value = 5399.385
$ws.Range("I85").Value = 1064.1428
$ws.Range("J85").Value = 10457.167
$ws.Range("K85").Value = 1064.1428
$ws.Range("L85").Value = 10457.167
$ws.Range("M85").Value = 183.8571999999999
$ws.Range("N85").Value = -12953.167

# row 93
$ws.Range("H93").Value = 1585.5
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# row 132
$ws.Range("H132").Value = 2558.6
$ws.Range("I132").Value = 2443.25
$ws.Range("J132").Value = 2635.5
$ws.Range("K132").Value = 7329.75
$ws.Range("L132").Value = 7906.5
$ws.Range("M132").Value = -4799.75
$ws.Range("N132").Value = -12966.5

# row 136
$ws.Range("H136").Value = 3649.6667
$ws.Range("J136").Value = 4349.375
$ws.Range("L136").Value = 13048.125
$ws.Range("N136").Value = -18148.125

$ws = $wb.Worksheets.Item("WVR")
# row 34
$ws.Range("H34").Value = 40008
$ws.Range("I34").Value = 40008
$ws.Range("K34").Value = 40008
$ws.Range("M34").Value = -39805

# row 37
$ws.Range("H37").Value = 25026
$ws.Range("I37").Value = 25026
$ws.Range("K37").Value = 25026
$ws.Range("M37").Value = -24823

# row 41
$ws.Range("H41").Value = 11131.625
$ws.Range("I41").Value = 10998
$ws.Range("J41").Value = 11150.714
$ws.Range("K41").Value = 10998
$ws.Range("L41").Value = 11150.714
$ws.Range("M41").Value = -10608
$ws.Range("N41").Value = -11930.714

# row 42
$ws.Range("H42").Value = 20000
$ws.Range("I42").Value = 20000
$ws.Range("K42").Value = 20000
$ws.Range("M42").Value = -19622

# row 100
$ws.Range("H100").Value = 891.6
$ws.Range("I100").Value = 413.14285
$ws.Range("K100").Value = 826.2857
$ws.Range("M100").Value = -285.2857

# row 136
$ws.Range("H136").Value = 48233.78
$ws.Range("J136").Value = 6522.3
$ws.Range("L136").Value = 19566.9
$ws.Range("N136").Value = -24666.9

